$wb = $excel.ActiveWorkbook

# Overview sheet: G4 -> Latest HO Xliff Generate Date for 9c23dcd6 row
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-08-25 18:47:09"

# zh-cn sheet: H4 -> Correspond Handoff Datetime, K4 -> Correspond Handback DateTime
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H4").Value = "2016-08-25 18:46:58"
$wsZhCn.Range("K4").Value = "2016-08-25 18:47:39"

# de-de sheet: K4 -> Correspond Handback DateTime (H4 shares the same shared string as Overview G4,
# so it updates automatically as a side effect in the diff but we set it explicitly too for clarity)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H4").Value = "2016-08-25 18:47:09"
$wsDeDe.Range("K4").Value = "2016-08-25 18:47:46"
